# Update the LR-pairs TPM-derived statistics (columns G..T, rows 2..16)
# on the active worksheet to reflect re-computed values using the new TPM data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Inflammatory-Mac)
$ws.Cells.Item(2, 7).Value = 197.5433703333333
$ws.Cells.Item(2, 8).Value = 592.6301109999999
$ws.Cells.Item(2, 9).Value = 0.3388703761585983
$ws.Cells.Item(2, 10).Value = 0.3388703761585982
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.011299
$ws.Cells.Item(2, 14).Value = 0.033897
$ws.Cells.Item(2, 15).Value = 0.3524439315012944
$ws.Cells.Item(2, 16).Value = 0.3524439315012944
$ws.Cells.Item(2, 17).Value = 2.232042541396333
$ws.Cells.Item(2, 18).Value = 20.088382872567
$ws.Cells.Item(2, 19).Value = 0.1194328076426589
$ws.Cells.Item(2, 20).Value = 0.1194328076426589

# Row 3 (ECs -> MuSCs)
$ws.Cells.Item(3, 7).Value = 197.5433703333333
$ws.Cells.Item(3, 8).Value = 592.6301109999999
$ws.Cells.Item(3, 9).Value = 0.3388703761585983
$ws.Cells.Item(3, 10).Value = 0.3388703761585982
$ws.Cells.Item(3, 15).Value = 0.01059504871227008
$ws.Cells.Item(3, 16).Value = 0.01059504871227008
$ws.Cells.Item(3, 17).Value = 0.06709889812322221
$ws.Cells.Item(3, 18).Value = 0.6038900831089999
$ws.Cells.Item(3, 19).Value = 0.003590348142545636
$ws.Cells.Item(3, 20).Value = 0.003590348142545635

# Row 4 (ECs -> Resolving-Mac)
$ws.Cells.Item(4, 7).Value = 197.5433703333333
$ws.Cells.Item(4, 8).Value = 592.6301109999999
$ws.Cells.Item(4, 9).Value = 0.3388703761585983
$ws.Cells.Item(4, 10).Value = 0.3388703761585982
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.02042033333333334
$ws.Cells.Item(4, 14).Value = 0.061261
$ws.Cells.Item(4, 15).Value = 0.6369610197864354
$ws.Cells.Item(4, 16).Value = 0.6369610197864355
$ws.Cells.Item(4, 17).Value = 4.033901469996778
$ws.Cells.Item(4, 18).Value = 36.305113229971
$ws.Cells.Item(4, 19).Value = 0.2158472203733937
$ws.Cells.Item(4, 20).Value = 0.2158472203733937

# Row 5 (FAPs -> Inflammatory-Mac)
$ws.Cells.Item(5, 9).Value = 0.1369374790620155
$ws.Cells.Item(5, 10).Value = 0.1369374790620154
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.011299
$ws.Cells.Item(5, 14).Value = 0.033897
$ws.Cells.Item(5, 15).Value = 0.3524439315012944
$ws.Cells.Item(5, 16).Value = 0.3524439315012944
$ws.Cells.Item(5, 17).Value = 0.9019681278806665
$ws.Cells.Item(5, 18).Value = 8.117713150925999
$ws.Cells.Item(5, 19).Value = 0.04826278349049291
$ws.Cells.Item(5, 20).Value = 0.04826278349049291

# Row 6 (FAPs -> MuSCs)
$ws.Cells.Item(6, 9).Value = 0.1369374790620155
$ws.Cells.Item(6, 10).Value = 0.1369374790620154
$ws.Cells.Item(6, 15).Value = 0.01059504871227008
$ws.Cells.Item(6, 16).Value = 0.01059504871227008
$ws.Cells.Item(6, 19).Value = 0.001450859261197519
$ws.Cells.Item(6, 20).Value = 0.001450859261197518

# Row 7 (FAPs -> Resolving-Mac)
$ws.Cells.Item(7, 9).Value = 0.1369374790620155
$ws.Cells.Item(7, 10).Value = 0.1369374790620154
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.02042033333333334
$ws.Cells.Item(7, 14).Value = 0.061261
$ws.Cells.Item(7, 15).Value = 0.6369610197864354
$ws.Cells.Item(7, 16).Value = 0.6369610197864355
$ws.Cells.Item(7, 17).Value = 1.630099108537556
$ws.Cells.Item(7, 18).Value = 14.670891976838
$ws.Cells.Item(7, 19).Value = 0.08722383631032501
$ws.Cells.Item(7, 20).Value = 0.08722383631032501

# Row 8 (Inflammatory-Mac -> Inflammatory-Mac)
$ws.Cells.Item(8, 7).Value = 148.824417
$ws.Cells.Item(8, 8).Value = 446.473251
$ws.Cells.Item(8, 9).Value = 0.2552967790580629
$ws.Cells.Item(8, 10).Value = 0.2552967790580629
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.011299
$ws.Cells.Item(8, 14).Value = 0.033897
$ws.Cells.Item(8, 15).Value = 0.3524439315012944
$ws.Cells.Item(8, 16).Value = 0.3524439315012944
$ws.Cells.Item(8, 17).Value = 1.681567087683
$ws.Cells.Item(8, 18).Value = 15.134103789147
$ws.Cells.Item(8, 19).Value = 0.08997780051084099
$ws.Cells.Item(8, 20).Value = 0.089977800510841

# Row 9 (Inflammatory-Mac -> MuSCs)
$ws.Cells.Item(9, 7).Value = 148.824417
$ws.Cells.Item(9, 8).Value = 446.473251
$ws.Cells.Item(9, 9).Value = 0.2552967790580629
$ws.Cells.Item(9, 10).Value = 0.2552967790580629
$ws.Cells.Item(9, 15).Value = 0.01059504871227008
$ws.Cells.Item(9, 16).Value = 0.01059504871227008
$ws.Cells.Item(9, 17).Value = 0.050550693641
$ws.Cells.Item(9, 18).Value = 0.454956242769
$ws.Cells.Item(9, 19).Value = 0.002704881810205829
$ws.Cells.Item(9, 20).Value = 0.002704881810205829

# Row 10 (Inflammatory-Mac -> Resolving-Mac)
$ws.Cells.Item(10, 7).Value = 148.824417
$ws.Cells.Item(10, 8).Value = 446.473251
$ws.Cells.Item(10, 9).Value = 0.2552967790580629
$ws.Cells.Item(10, 10).Value = 0.2552967790580629
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.02042033333333334
$ws.Cells.Item(10, 14).Value = 0.061261
$ws.Cells.Item(10, 15).Value = 0.6369610197864354
$ws.Cells.Item(10, 16).Value = 0.6369610197864355
$ws.Cells.Item(10, 17).Value = 3.039044203279
$ws.Cells.Item(10, 18).Value = 27.351397829511
$ws.Cells.Item(10, 19).Value = 0.162614096737016
$ws.Cells.Item(10, 20).Value = 0.162614096737016

# Row 11 (MuSCs -> Inflammatory-Mac)
$ws.Cells.Item(11, 7).Value = 35.426853
$ws.Cells.Item(11, 8).Value = 106.280559
$ws.Cells.Item(11, 9).Value = 0.06077202683121193
$ws.Cells.Item(11, 10).Value = 0.06077202683121192
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.011299
$ws.Cells.Item(11, 14).Value = 0.033897
$ws.Cells.Item(11, 15).Value = 0.3524439315012944
$ws.Cells.Item(11, 16).Value = 0.3524439315012944
$ws.Cells.Item(11, 17).Value = 0.4002880120469999
$ws.Cells.Item(11, 18).Value = 3.602592108422999
$ws.Cells.Item(11, 19).Value = 0.02141873206169448
$ws.Cells.Item(11, 20).Value = 0.02141873206169448

# Row 12 (MuSCs -> MuSCs)
$ws.Cells.Item(12, 7).Value = 35.426853
$ws.Cells.Item(12, 8).Value = 106.280559
$ws.Cells.Item(12, 9).Value = 0.06077202683121193
$ws.Cells.Item(12, 10).Value = 0.06077202683121192
$ws.Cells.Item(12, 15).Value = 0.01059504871227008
$ws.Cells.Item(12, 16).Value = 0.01059504871227008
$ws.Cells.Item(12, 17).Value = 0.012033321069
$ws.Cells.Item(12, 18).Value = 0.108299889621
$ws.Cells.Item(12, 19).Value = 0.0006438825846200749
$ws.Cells.Item(12, 20).Value = 0.0006438825846200748

# Row 13 (MuSCs -> Resolving-Mac)
$ws.Cells.Item(13, 7).Value = 35.426853
$ws.Cells.Item(13, 8).Value = 106.280559
$ws.Cells.Item(13, 9).Value = 0.06077202683121193
$ws.Cells.Item(13, 10).Value = 0.06077202683121192
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.02042033333333334
$ws.Cells.Item(13, 14).Value = 0.061261
$ws.Cells.Item(13, 15).Value = 0.6369610197864354
$ws.Cells.Item(13, 16).Value = 0.6369610197864355
$ws.Cells.Item(13, 17).Value = 0.7234281472110001
$ws.Cells.Item(13, 18).Value = 6.510853324899
$ws.Cells.Item(13, 19).Value = 0.03870941218489736
$ws.Cells.Item(13, 20).Value = 0.03870941218489736

# Row 14 (Resolving-Mac -> Inflammatory-Mac)
$ws.Cells.Item(14, 7).Value = 121.3248153333333
$ws.Cells.Item(14, 8).Value = 363.974446
$ws.Cells.Item(14, 9).Value = 0.2081233388901116
$ws.Cells.Item(14, 10).Value = 0.2081233388901115
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.011299
$ws.Cells.Item(14, 14).Value = 0.033897
$ws.Cells.Item(14, 15).Value = 0.3524439315012944
$ws.Cells.Item(14, 16).Value = 0.3524439315012944
$ws.Cells.Item(14, 17).Value = 1.370849088451333
$ws.Cells.Item(14, 18).Value = 12.337641796062
$ws.Cells.Item(14, 19).Value = 0.07335180779560715
$ws.Cells.Item(14, 20).Value = 0.07335180779560715

# Row 15 (Resolving-Mac -> MuSCs)
$ws.Cells.Item(15, 7).Value = 121.3248153333333
$ws.Cells.Item(15, 8).Value = 363.974446
$ws.Cells.Item(15, 9).Value = 0.2081233388901116
$ws.Cells.Item(15, 10).Value = 0.2081233388901115
$ws.Cells.Item(15, 15).Value = 0.01059504871227008
$ws.Cells.Item(15, 16).Value = 0.01059504871227008
$ws.Cells.Item(15, 17).Value = 0.04120999560822222
$ws.Cells.Item(15, 18).Value = 0.370889960474
$ws.Cells.Item(15, 19).Value = 0.002205076913701027
$ws.Cells.Item(15, 20).Value = 0.002205076913701027

# Row 16 (Resolving-Mac -> Resolving-Mac)
$ws.Cells.Item(16, 7).Value = 121.3248153333333
$ws.Cells.Item(16, 8).Value = 363.974446
$ws.Cells.Item(16, 9).Value = 0.2081233388901116
$ws.Cells.Item(16, 10).Value = 0.2081233388901115
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.02042033333333334
$ws.Cells.Item(16, 14).Value = 0.061261
$ws.Cells.Item(16, 15).Value = 0.6369610197864354
$ws.Cells.Item(16, 16).Value = 0.6369610197864355
$ws.Cells.Item(16, 17).Value = 2.477493170711778
$ws.Cells.Item(16, 18).Value = 22.297438536406
$ws.Cells.Item(16, 19).Value = 0.1325664541808033
$ws.Cells.Item(16, 20).Value = 0.1325664541808033
